$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the completed "stock-check" (kiem kho) flow
$ws.Name = "Template kiểm kho"

# Move/restore the active selection to D26, matching the saved view state
$ws.Range("D26").Select() | Out-Null
